# Add team record (Wins / Losses / Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column headers in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, centered,
# thin border) by copying the format from the existing AC1 header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows (2-55): same record for every player/row in the sheet.
$ws.Range("AD2:AD55").Value = 61
$ws.Range("AE2:AE55").Value = 101
$ws.Range("AF2:AF55").Value = 0
